$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the row containing account 005890232 / TAYLA / 26377.15 (row 3)
$ws.Rows.Item(3).Delete()
